# Daily attendance processing - 2026-02-14 16:35:16 UTC
# Move "Administrator" to the front of the "Recorded By" (column G) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains("Administrator")) {
        $parts = $val -split ", " | ForEach-Object { $_.Trim() }
        $others = $parts | Where-Object { $_ -ne "Administrator" }
        $newParts = @("Administrator") + $others
        $newVal = [string]::Join(", ", $newParts)
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}
